$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Travail")

# Copy the date-formatted style from row 9 onto the date cells of rows 10-11
# so the new cells keep the existing date number format (style index 7)
# instead of Excel inventing a brand-new custom format.
$ws.Range("B9").Copy()
$ws.Range("B10:B11").PasteSpecial(-4122)
$ws.Range("F9").Copy()
$ws.Range("F10:F11").PasteSpecial(-4122)

# Row 10
$ws.Range("B10").Value = 45384
$ws.Range("C10").Value = "2h00"
$ws.Range("D10").Value = "Programmation de la rom et création du test bench"
$ws.Range("F10").Value = 45384
$ws.Range("G10").Value = "2h00"
$ws.Range("H10").Value = "Programmation de la rom et création du test bench"

# Row 11
$ws.Range("B11").Value = 45391
$ws.Range("C11").Value = "1h30"
$ws.Range("D11").Value = "Debug du code grâce à la carte et correction de certain bug"
$ws.Range("F11").Value = 45391
$ws.Range("G11").Value = "1h30"
$ws.Range("H11").Value = "Debug du code grâce à la carte et correction de certain bug"

$ws.Range("J12").Select()
